# Update "想去人数" (interested-count) figures on the gh-pages data refresh.
# Values below mirror the diff for the regenerated 苏州-漫展信息.xlsx output.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 1211
$ws1.Range("F4").Value  = 295
$ws1.Range("F6").Value  = 17
$ws1.Range("F7").Value  = 12420
$ws1.Range("F8").Value  = 65
$ws1.Range("F11").Value = 3
$ws1.Range("F12").Value = 178
$ws1.Range("F13").Value = 12270
$ws1.Range("F14").Value = 4861
$ws1.Range("F15").Value = 4749
$ws1.Range("F16").Value = 146
$ws1.Range("F17").Value = 67
$ws1.Range("F20").Value = 956
$ws1.Range("F23").Value = 171

# --- Sheet "演出" ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 7

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 1211
$ws4.Range("F4").Value  = 295
$ws4.Range("F6").Value  = 7
$ws4.Range("F8").Value  = 17
$ws4.Range("F9").Value  = 12420
$ws4.Range("F10").Value = 65
$ws4.Range("F13").Value = 3
$ws4.Range("F14").Value = 179
$ws4.Range("F15").Value = 12270
$ws4.Range("F16").Value = 4861
$ws4.Range("F17").Value = 4749
$ws4.Range("F18").Value = 146
$ws4.Range("F19").Value = 67
$ws4.Range("F22").Value = 956
$ws4.Range("F25").Value = 171
